$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Append a new test-log row (row 14) -----------------------------------
# Copy the formatting of the most similar existing rows/cells first (so the
# new cells reuse the same cell styles as their neighbours instead of Excel
# minting brand-new style records), then fill in the actual content.

$ws.Range("A13").Copy()
$ws.Range("A14").PasteSpecial(-4122)   # xlPasteFormats

$ws.Range("B13").Copy()
$ws.Range("B14").PasteSpecial(-4122)

$ws.Range("C13").Copy()
$ws.Range("C14").PasteSpecial(-4122)

$ws.Range("D13").Copy()
$ws.Range("D14").PasteSpecial(-4122)

$ws.Range("E13").Copy()
$ws.Range("E14").PasteSpecial(-4122)

$ws.Range("F8").Copy()
$ws.Range("F14").PasteSpecial(-4122)

# Date (2017-09-13 -> Excel serial 42991, same representation the sheet
# already uses for every other date in column A).
$ws.Range("A14").Value = 42991

# Fill content in the same left-to-right-ish order the original author used
# (this also controls the order new entries land in sharedStrings.xml).
$ws.Range("F14").Value = "Wachen stehen beim erreichten Gegenstand ineinander; Wachen drehen sich noch nicht korrekt zur Wand; Wachen drehen sich noch nicht korrekt zum Spieler; manchmal Absturz, weil sich eine Coroutine aufhängt…? (Wachen patrouillieren & Gegenstand rotieren)"
$ws.Range("E14").Value = "Wachen können schießen & Spieler jagen"
$ws.Range("C14").Value = "Entwicklung"
$ws.Range("D14").Value = "Anna Franziska"
$ws.Range("B14").Value = "DiscordiaAgency_Demo_2017_09_13.exe"

# Row height for the new (wrapped, multi-line) row.
$ws.Rows.Item(14).RowHeight = 120

# Move the selection the way the author's Excel session ended up: cell B15
# (just past the freshly typed row).
$ws.Range("B15").Select()
